# Update crypto price / 1h-volume-change figures (GitHub Actions data refresh).
# NumberFormat is forced to Text ("@") before assignment so values such as
# "212.48" or "0.900" are kept as literal strings (matching the source
# workbook's inlineStr cells) instead of being auto-converted to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.541.92'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.639.21'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.48'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.536'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.55%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.02'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.256'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0891'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.79'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.645.86'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.02'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.05'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.561.00'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.89'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.47%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.81%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.88'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.99%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.29'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.62%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.60'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.426.27'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.55%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.97%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.876'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.81%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.900'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +15.12%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.75%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.26'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.92'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.780.95'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.22'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.48%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.65%  '
